{"js": "// The document's final paragraph (\"Motif du refoulement : \u00ab=Motif\u00bb\" \u2014 a\n// MERGEFIELD block) is removed in its entirety as part of refactoring how\n// the repetitive publipostage (mail-merge) blocks are read.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph to remove is the last one in the body.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.delete();\n\nawait context.sync();\n", "ps1": "# The document's final paragraph (\"Motif du refoulement : \u00ab=Motif\u00bb\" \u2014 a\n# MERGEFIELD block) is removed in its entirety as part of refactoring how\n# the repetitive publipostage (mail-merge) blocks are read.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.Delete()\n"}
